# CardStats.xlsx update — refreshed card-count figures (Chorus, various
# aurics, archangel shield) on the "Energy Distribution" / "Rarity
# Distribution" / "Type Distribution" tables on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Energy Distribution table (column E, "3" energy-cost row group) ---
$ws.Range("E8").Value  = 9   # was 11
$ws.Range("E9").Value  = 31  # was 24
$ws.Range("E10").Value = 8   # was 11
$ws.Range("E11").Value = 2   # was 5
$ws.Range("E12").Value = 4   # was 6
# E13 = SUM(E7:E12) recalculates automatically (59 -> 56)

# --- Rarity Distribution table (column K, "Hierophant" column) ---
$ws.Range("K8").Value = 14  # was 18
$ws.Range("K9").Value = 27  # was 29
# K11 = SUM(K7:K10) recalculates automatically (63 -> 57)

# --- Type Distribution table (column K, "Hierophant" column) ---
$ws.Range("K17").Value = 18  # was 20
$ws.Range("K18").Value = 34  # was 29
$ws.Range("K19").Value = 5   # was 10
# K20 = SUM(K17:K19) recalculates automatically (59 -> 57)

# Move the cursor/selection to where the author ended up (K19) and nudge
# the window size slightly, matching the end-of-session UI state.
$ws.Range("K19").Select()
$excel.ActiveWindow.Height = 371.5
